$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (ProjektID) values change from numeric to text codes
$ws.Range("A4").Value = "A1"
$ws.Range("A5").Value = "A2"

# Update the active selection to A6
$ws.Range("A6").Select()
